# Update Work Week and Social Spending
# - Data sheet: refresh existing "GDP per Capita" values in column E (rows 2-192)
# - Data sheet: append new rows for years 2011-2016

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Map of row number -> new Data value (column E). Values are plain numeric-looking
# text, so each cell is re-written as a *text* cell (matching the source data, which
# stores every "Data" entry as a shared string, not a number).
$dataUpdates = [ordered]@{
    2 = "840"
    52 = "805"
    95 = "1159"
    132 = "1049"
    133 = "1078"
    134 = "1106"
    135 = "1135"
    136 = "1167"
    137 = "1195"
    138 = "1218"
    139 = "1235"
    140 = "1251"
    141 = "1262"
    142 = "1274"
    143 = "1294"
    144 = "1411"
    145 = "1406"
    146 = "1427"
    147 = "1398"
    148 = "1369"
    149 = "1165"
    150 = "1114"
    151 = "1178"
    152 = "1172"
    153 = "1202"
    154 = "1278"
    155 = "1333"
    156 = "1248"
    157 = "1132"
    158 = "1290"
    159 = "1304"
    160 = "1285"
    161 = "1267"
    162 = "1207"
    163 = "1223"
    164 = "1293"
    165 = "1333"
    166 = "1420"
    167 = "1471"
    168 = "1481"
    169 = "1500"
    170 = "1556"
    171 = "1588"
    172 = "1634"
    173 = "1688.15989504014"
    174 = "1792.67077363815"
    175 = "1894.85862029566"
    176 = "2018.35070301574"
    177 = "2167.04204895783"
    178 = "2323.04638353782"
    179 = "2466.87354952478"
    180 = "2564.28445564975"
    181 = "2641.3484975664"
    182 = "2773.10163510792"
    183 = "2915.68614024151"
    184 = "3070.6934784829"
    185 = "3242.2854696511"
    186 = "3437.6644223215"
    187 = "3637.65637930508"
    188 = "3828.90767342617"
    189 = "4037.23285915626"
    190 = "4200.23523093943"
    191 = "4360.06172503739"
    192 = "4571.84858405268"
}

foreach ($row in $dataUpdates.Keys) {
    $cell = $ws.Range("E$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dataUpdates[$row]
    $cell.ClearFormats()
}

# Append the new 2011-2016 rows at the bottom of the table.
$newRows = @(
    @{ Row = 193; Year = 2011; Value = "4786" }
    @{ Row = 194; Year = 2012; Value = "4984" }
    @{ Row = 195; Year = 2013; Value = "5200" }
    @{ Row = 196; Year = 2014; Value = "5455" }
    @{ Row = 197; Year = 2015; Value = "5763" }
    @{ Row = 198; Year = 2016; Value = "6062" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = 704
    $ws.Cells.Item($r, 2).Value = "Vietnam"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $entry.Year
    $cell = $ws.Range("E$r")
    $cell.NumberFormat = "@"
    $cell.Value = $entry.Value
    $cell.ClearFormats()
}
